$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.701.31'
$ws.Range("E2").Value = '  +0.27%  '

$ws.Range("D3").Value = '3.104.27'
$ws.Range("E3").Value = '  +1.10%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '523.17'
$ws.Range("E5").Value = '  +1.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.50'
$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").Value = '3.108.39'
$ws.Range("E8").Value = '  +1.31%  '

$ws.Range("E9").Value = '  +0.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.27'
$ws.Range("E10").Value = '  -0.20%  '

$ws.Range("E11").Value = '  +0.78%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.382'
$ws.Range("E12").Value = '  +2.39%  '

$ws.Range("D13").Value = '3.643.50'
$ws.Range("E13").Value = '  +1.11%  '

$ws.Range("E14").Value = '  +1.68%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.15'
$ws.Range("E15").Value = '  +2.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000164'
$ws.Range("E16").Value = '  +0.70%  '

$ws.Range("D17").Value = '57.826.29'

$ws.Range("D18").Value = '3.107.18'
$ws.Range("E18").Value = '  +1.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.10'
$ws.Range("E19").Value = '  +0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.80'
$ws.Range("E20").Value = '  -1.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.04'
$ws.Range("E21").Value = '  -0.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '336.87'
$ws.Range("E22").Value = '  +0.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.25%  '

$ws.Range("E24").Value = '  +1.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.54'
$ws.Range("E25").Value = '  +0.85%  '

$ws.Range("E26").Value = '  -1.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.04%  '

$ws.Range("D28").Value = '0.0₃0927'
$ws.Range("E28").Value = '  +1.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.55'
$ws.Range("E29").Value = '  +2.89%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.22'
$ws.Range("E31").Value = '  +0.34%  '

$ws.Range("E32").Value = '  +2.09%  '

$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.89'
$ws.Range("E33").Value = '  +0.32%  '

$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.20'
$ws.Range("E34").Value = '  +2.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '154.20'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.61'
$ws.Range("E36").Value = '  +3.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.08'
$ws.Range("E37").Value = '  +2.99%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.98'
$ws.Range("E38").Value = '  -0.27%  '

$ws.Range("E39").Value = '  +1.96%  '

$ws.Range("E40").Value = '  -1.15%  '

$ws.Range("D41").Value = '3.152.78'
$ws.Range("E41").Value = '  +1.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.685'
$ws.Range("E42").Value = '  +4.38%  '

$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '36.94'
$ws.Range("E43").Value = '  -0.03%  '

$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.88'
$ws.Range("E44").Value = '  -0.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("E46").Value = '  +5.53%  '

$ws.Range("D47").Value = '2.281.13'
$ws.Range("E47").Value = '  +0.55%  '

$ws.Range("E48").Value = '  +0.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.982'
$ws.Range("E49").Value = '  +5.81%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.63'
$ws.Range("E50").Value = '  +1.16%  '

$ws.Range("E51").Value = '  +2.16%  '
